$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 67-70: columns F:V get cyclically rotated "up" (67<-68, 68<-69,
#    69<-70, 70<-67). Columns A:E (id / league / season / date) stay put for
#    each row. Read all four F:V rows into memory first, then write them back
#    rotated, so nothing gets clobbered mid-update.
# ---------------------------------------------------------------------------
$row67 = $ws.Range("F67:V67").Value2
$row68 = $ws.Range("F68:V68").Value2
$row69 = $ws.Range("F69:V69").Value2
$row70 = $ws.Range("F70:V70").Value2

$ws.Range("F67:V67").Value2 = $row68
$ws.Range("F68:V68").Value2 = $row69
$ws.Range("F69:V69").Value2 = $row70
$ws.Range("F70:V70").Value2 = $row67

# ---------------------------------------------------------------------------
# 2) Append new row 91 (match fixture "0" which is A=90) with its data, and
#    copy number/cell formatting from the row above (90) so styles match.
# ---------------------------------------------------------------------------
$ws.Range("A90:V90").Copy()
$ws.Range("A91:V91").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A91").Value = 90
$ws.Range("B91").Value = "denmark"
$ws.Range("C91").Value = "2nd-division"
$ws.Range("D91").Value = "2023-2024"
$ws.Range("E91").Value = 45247.75
$ws.Range("F91").Value = "F. Amager"
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = "Nykobing"
$ws.Range("I91").Value = 4
$ws.Range("J91").Value = 2.51
$ws.Range("K91").Value = "16/11/2023 06:12"
$ws.Range("L91").Value = 2.59
$ws.Range("M91").Value = "17/11/2023 16:42"
$ws.Range("N91").Value = 3.35
$ws.Range("O91").Value = "16/11/2023 06:12"
$ws.Range("P91").Value = 3.48
$ws.Range("Q91").Value = "17/11/2023 16:42"
$ws.Range("R91").Value = 2.45
$ws.Range("S91").Value = "16/11/2023 06:12"
$ws.Range("T91").Value = 2.52
$ws.Range("U91").Value = "17/11/2023 16:42"
$ws.Range("V91").Value = "https://www.betexplorer.com/football/denmark/2nd-division/fremad-amager-nykobing/nTvp7xmi/"
